# "Adicionados balanços concatenados em uma única planilha."
#
# This sheet holds several balance sheets that were concatenated into one
# worksheet. A handful of trailing quarter columns that were stubbed out
# with literal 0s when the sheets were stitched together are reset to
# empty-text placeholders (matching the blank cells used everywhere else
# in this concatenated table), and one date serial (AF57) gets a small
# fractional correction.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Small correction to a date serial value that picked up a fractional
# remainder during concatenation.
$ws.Range("AF57").Value = 44888.008

# Trailing placeholder columns (all zeros) that belong to quarters beyond
# the real data for these rows -- reset them to the same blank
# empty-text placeholder used throughout the rest of the concatenated
# table (rather than deleting the cells outright).
$blankRanges = @(
    "AZ57:BT57",
    "AZ58:BT58",
    "AE64:AZ64",
    "AZ71:BT71",
    "AZ72:BT72",
    "AZ73:BT73",
    "AZ77:BT77",
    "AZ78:BT78",
    "AE79:BT79"
)

foreach ($addr in $blankRanges) {
    $rng = $ws.Range($addr)
    # Entering a leading apostrophe forces text type with an empty
    # display value (the same "blank" placeholder already used by every
    # other empty cell in this sheet), then the style is reset so no
    # stray quote-prefix formatting is left behind.
    $rng.Value = "'"
    $rng.Style = "Normal"
}
